# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-suffix columns to "_FV2410" / "_FV2504"
# 2) Freeze the header row (top row) on the sheet
# 3) Turn the data range into an Excel Table ("Table1") with autofilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename headers ----------------------------------------------------
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = ($text -replace "_old$", "_FV2410")
}

$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = ($text -replace "_new$", "_FV2504")
}

# --- 2) Freeze the top row -------------------------------------------------
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- 3) Convert the used range into an Excel Table --------------------------
$tableRange = $ws.Range("A1:U81")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

$wb.Save()
